# BOT; UPDATE DATA
# Appends one new day of data to the "相談件数" sheet:
#   - the old row 106 (the "※4/8..." footnote row) is pushed down to row 107
#   - the freed-up row 106 gets a new data row (date 2020-05-10 / serial 43961)
#   - dimension / print area / selection are refreshed to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the existing row 106 (footnote row) down to row 107, inserting a
# fresh row 106 that copies formatting from the row above it (row 105),
# exactly like Excel does when a row is inserted above an existing row.
$ws.Rows.Item(106).Insert()

# Fill the newly inserted row 106 with the new day's figures.
$ws.Range("A106").Value = 43961
$ws.Range("B106").Value = 394
$ws.Range("C106").Value = 35779
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 7232

# Scroll the frozen window so row 83 is the first visible (scrolled) row,
# then restore the active selection to match the refreshed sheet (A106).
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("A106").Select()

# Extend the print area by one row to keep a blank buffer row below data
# (mirrors the existing convention: Print_Area bottom = dimension bottom + 1).
$name = $wb.Names.Item(1)
$name.RefersTo = "=相談件数!`$A`$1:`$E`$108"
